$d = $word.ActiveDocument

# --- Locate the target paragraph: the empty paragraph that currently
# holds the "_GoBack" bookmark, right after the "Tricky Points:" heading
# and right before the "Overlapping Sub - Question:" paragraph. ---
$bm = $d.Bookmarks.Item("_GoBack")
$bmPara = $bm.Range.Paragraphs.Item(1)
$bmIndex = $bmPara.Index
$insertPoint = $d.Range($bmPara.Range.Start, $bmPara.Range.Start)

# Remember where the "Overlapping Sub - Question:" paragraph (which
# follows the bookmark paragraph) is, so we can relocate the bookmark
# into the empty paragraph that follows it once the new runs exist.
# (Indices are stable here: the edit only adds runs of text, not
# paragraphs, so paragraph numbering doesn't shift.)
$afterOverlapIndex = $bmIndex + 2

# --- Insert the three new runs of text into the (currently empty)
# bookmark paragraph, preserving the exact run-level formatting from
# the target revision (plain / highlighted-term / plain). ---
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">We need to figure out whether the problem has the Best Sub - Structure when trying to use Dynamic Programming Algorithm. Normally, </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="C00000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:shd w:val="clear" w:color="FFFFFF" w:fill="D9D9D9"/></w:rPr><w:t>when solving two problems and they would not share the same resource, then we can call two problems independent.</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insertPoint.InsertXML($xml)

# --- Move the "_GoBack" bookmark: it no longer marks the (now filled
# in) paragraph; instead it marks the empty paragraph that follows the
# "Overlapping Sub - Question:" paragraph. ---
$bm.Delete()
$afterOverlapPara = $d.Paragraphs.Item($afterOverlapIndex)
$d.Bookmarks.Add("_GoBack", $afterOverlapPara.Range)
